# Updated cryptos list on Wed Feb  7 17:14:41 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.267.27"
$ws.Range("E2").Value = "  +0.18%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.412.74"
$ws.Range("E3").Value = "  +2.57%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'306.97"
$ws.Range("E5").Value = "  +1.38%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'97.16"
$ws.Range("E6").Value = "  +1.69%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.506"
$ws.Range("E7").Value = "  +0.28%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.487"
$ws.Range("E9").Value = "  -2.20%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "'35.04"
$ws.Range("E10").Value = "  +2.60%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.126"
$ws.Range("E11").Value = "  +3.37%  "

# Row 12 - Dogecoin
$ws.Range("D12").Value = "'0.0795"
$ws.Range("E12").Value = "  +0.94%  "

# Row 13 - Chainlink
$ws.Range("D13").Value = "'18.63"
$ws.Range("E13").Value = "  +0.39%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'6.87"
$ws.Range("E14").Value = "  +1.62%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.768.82"
$ws.Range("E15").Value = "  +1.98%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.404.53"
$ws.Range("E16").Value = "  +2.56%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "'0.825"
$ws.Range("E17").Value = "  +3.47%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "43.306.90"
$ws.Range("E18").Value = "  +0.35%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "'6.41"
$ws.Range("E19").Value = "  +2.87%  "

# Row 20 - InternetComputer(DFINITY)
$ws.Range("D20").Value = "'12.10"
$ws.Range("E20").Value = "  -0.96%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0₃0901"
$ws.Range("E21").Value = "  +1.14%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "'68.42"
$ws.Range("E22").Value = "  +0.36%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "'237.69"
$ws.Range("E23").Value = "  +0.99%  "

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  +1.89%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "'2.47"
$ws.Range("E25").Value = "  +1.55%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.01%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'24.82"
$ws.Range("E27").Value = "  +0.90%  "

# Row 28 - Toncoin
$ws.Range("D28").Value = "'2.34"
$ws.Range("E28").Value = "  -0.82%  "

# Row 29 - Cosmos
$ws.Range("D29").Value = "'9.41"
$ws.Range("E29").Value = "  +2.65%  "

# Row 30 - InjectiveProtocol
$ws.Range("D30").Value = "'32.30"
$ws.Range("E30").Value = "  +3.31%  "

# Row 31 - FirstDigitalUSD -> Filecoin
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'5.15"
$ws.Range("E31").Value = "  +2.74%  "

# Row 32 - Filecoin -> Celestia
$ws.Range("B32").Value = "Celestia"
$ws.Range("C32").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D32").Value = "'18.44"
$ws.Range("E32").Value = "  +7.18%  "

# Row 33 - Kaspa -> FirstDigitalUSD
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.09%  "

# Row 34 - Celestia -> Kaspa
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "'0.114"
$ws.Range("E34").Value = "  +13.71%  "

# Row 35 - Hedera
$ws.Range("D35").Value = "'0.0746"
$ws.Range("E35").Value = "  +3.12%  "

# Row 36 - Monero -> LidoDAOToken
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'3.04"
$ws.Range("E36").Value = "  +10.52%  "

# Row 37 - ARBITRUM -> Monero
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "'131.32"
$ws.Range("E37").Value = "  +16.59%  "

# Row 38 - LidoDAOToken -> ARBITRUM
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'1.88"
$ws.Range("E38").Value = "  +1.77%  "

# Row 39 - RenderToken
$ws.Range("D39").Value = "'4.37"
$ws.Range("E39").Value = "  +0.15%  "

# Row 40 - WEMIXToken
$ws.Range("E40").Value = "  -1.62%  "

# Row 41 - Stellar
$ws.Range("E41").Value = "  -0.55%  "

# Row 42 - EnergySwap
$ws.Range("D42").Value = "'21.19"
$ws.Range("E42").Value = "  -6.40%  "

# Row 43 - Maker
$ws.Range("D43").Value = "1.947.28"
$ws.Range("E43").Value = "  +0.31%  "

# Row 44 - VeChain
$ws.Range("D44").Value = "'0.0282"
$ws.Range("E44").Value = "  +0.29%  "

# Row 45 - ApeXProtocol
$ws.Range("D45").Value = "'2.15"
$ws.Range("E45").Value = "  +1.20%  "

# Row 46 - NEARProtocol
$ws.Range("D46").Value = "'2.80"
$ws.Range("E46").Value = "  +2.39%  "

# Row 47 - FraxShare
$ws.Range("D47").Value = "'9.35"
$ws.Range("E47").Value = "  -4.54%  "

# Row 48 - RocketPoolETH
$ws.Range("D48").Value = "2.618.84"
$ws.Range("E48").Value = "  +1.49%  "

# Row 49 - Stacks
$ws.Range("D49").Value = "'1.54"
$ws.Range("E49").Value = "  +2.45%  "

# Row 50 - MultiversX
$ws.Range("D50").Value = "'52.44"
$ws.Range("E50").Value = "  -0.93%  "

# Row 51 - BitcoinSV
$ws.Range("D51").Value = "'72.03"
$ws.Range("E51").Value = "  -0.27%  "

